$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.759.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.005.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.03%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.98%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.96%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.008.06"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.362"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.548.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.125"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.783.39"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.028.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.05%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000148"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.99%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "391.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.90"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.65"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.81%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.06%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.462"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.186"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0963"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.74"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.09"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.95%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.98"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.69%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.506.79"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.30%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.55"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.11%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.51"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.94%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.662"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.34%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0595"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.90%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.18%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0245"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.13"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.37%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.94"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.99%  "

# Row 49
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "271.38"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.56%  "

# Row 50
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.50"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.48%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0940"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.38%  "
